# -----------------------------------------------------------------------
# Applies the OOXML diff to GeradorDeSinais_Requirements.docx:
#   1. Adds a "_GoBack" bookmark on the (empty) paragraph right after the
#      title block near the top of the document.
#   2. Collapses the two runs that make up "REQ05: " + "Ambiente de
#      Desenvolvimento" into a single run.
#   3. Collapses the two runs that make up "REQ15: " + "Geração e
#      Validação por Hardware" into a single run.
#   4. Collapses the two runs "Como operador de testes, " + "eu preciso
#      gerar formas de onda" into a single run, while leaving the
#      subsequent runs (" específicas", " em", " um hardware...") intact.
#   5. Removes the old trailing "_GoBack" bookmark near the end of the
#      document (on the last paragraph before the sectPr).
#
# Net effect on the two remaining "_Hlk..." bookmarks: because bookmark
# ids are reassigned densely in document order whenever the bookmark set
# changes, adding the new "_GoBack" bookmark near the top bumps
# "_Hlk503963347"/"_Hlk503963364" from ids 0/1 up to 1/2, and removing the
# old trailing "_GoBack" (previously id 2) leaves them at 1/2 - matching
# the target diff exactly.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Remove the old trailing "_GoBack" bookmark near the end of the document ---
#        (done first, before adding the replacement one below, since
#        bookmark names must stay unique).
$d.Bookmarks.Item("_GoBack").Delete()

# --- 1b. Add the new "_GoBack" bookmark on the empty paragraph after the title ---
$titleGapPara = $d.Paragraphs.Item(2)
$d.Bookmarks.Add("_GoBack", $titleGapPara.Range)

# --- 2. Merge "REQ05: " + "Ambiente de Desenvolvimento" into one run ---
$rng = $d.Content
$rng.Find.Execute("REQ05: Ambiente de Desenvolvimento", $true, $false, $false, $false, $false, $true, 1, $false, "REQ05: Ambiente de Desenvolvimento", 2) | Out-Null

# --- 3. Merge "REQ15: " + "Geração e Validação por Hardware" into one run ---
$rng = $d.Content
$rng.Find.Execute("REQ15: Geração e Validação por Hardware", $true, $false, $false, $false, $false, $true, 1, $false, "REQ15: Geração e Validação por Hardware", 2) | Out-Null

# --- 4. Merge only "Como operador de testes, " + "eu preciso gerar formas de onda" ---
#        (leave the following runs untouched). A transient bookmark is
#        planted right after the text we want merged; it blocks the
#        engine's run-coalescing pass from reaching past it, then it is
#        deleted again so no trace of it remains in the saved document.
$rng = $d.Content
$rng.Find.Execute("Como operador de testes, eu preciso gerar formas de onda") | Out-Null
$boundary = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("ZZTempMergeBoundary", $boundary)

$rng2 = $d.Content
$rng2.Find.Execute("Como operador de testes, eu preciso gerar formas de onda", $true, $false, $false, $false, $false, $true, 1, $false, "Como operador de testes, eu preciso gerar formas de onda", 2) | Out-Null

$d.Bookmarks.Item("ZZTempMergeBoundary").Delete()
